$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from row 2 (an un-tall data row) down into row 4 first,
# so the new row reuses the existing style entries instead of creating new ones.
$ws.Range("A2:C2").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new row 4 data: date, hours, description
$ws.Range("A4").Value = 41949
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Clean up HTML and partially integrate 960 grid system"

# Update selection to mirror the new active cell
$ws.Range("C4").Select()

# Update workbook view window position to match the diff (window was
# moved/resized on screen between saves)
$excel.ActiveWindow.Left = -27700
$excel.ActiveWindow.Top = 2080
